# Modification in the cover
# Update the date line on the cover (title) slide from
# "Valencia, July 2022" to "26th July 2022".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The date text lives in the "CustomShape 12" shape (11th shape on the
# cover slide) together with the author name and "Ph D. Dissertation".
$shp = $s.Shapes.Item(11)

$tr = $shp.TextFrame.TextRange
$full = $tr.Text

$target = "Valencia, July 2022"
$idx = $full.IndexOf($target)

if ($idx -ge 0) {
    $sub = $tr.Characters($idx + 1, $target.Length)
    $sub.Text = "26th July 2022"
}
